$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3, 0, 1, 1, 0, 1, 2, 3, 1, 2, 0, 0, 1, 1, 3, 3)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
